$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.094.97"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").Value = "1.666.21"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'216.54"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'0.5096"
$ws.Range("E6").Value = "  +2.80%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.2633"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").Value = "'0.06401"
$ws.Range("E9").Value = "  +5.39%  "
$ws.Range("D10").Value = "'21.52"
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").Value = "'0.07407"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").Value = "1.672.24"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "'4.507"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "'0.5792"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").Value = "'0.000008535"
$ws.Range("E15").Value = "  +4.06%  "
$ws.Range("D16").Value = "'64.07"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "26.162.60"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "'4.921"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "'10.80"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "'189.44"
$ws.Range("E21").Value = "  +4.12%  "
$ws.Range("D22").Value = "'6.195"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'145.11"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").Value = "'7.607"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("E26").Value = "  +5.87%  "
$ws.Range("D27").Value = "'15.57"
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("D28").Value = "'0.06322"
$ws.Range("E28").Value = "  +13.32%  "
$ws.Range("D29").Value = "'1.296"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").Value = "'1.315"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'3.522"
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("D32").Value = "'3.500"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").Value = "'1.633"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").Value = "'1.013"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").Value = "'0.6072"
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("D36").Value = "'2.365"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'2.647"
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("D38").Value = "'6.161"
$ws.Range("E38").Value = "  +5.05%  "
$ws.Range("D39").Value = "'0.01604"
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("D40").Value = "1.074.47"
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").Value = "'0.8621"
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").Value = "'1.009"
$ws.Range("D43").Value = "'100.97"
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("D44").Value = "1.814.22"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("E45").Value = "  +8.48%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "'8.055"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.4295"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05200"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "'5.917"
$ws.Range("E51").Value = "  +6.36%  "
